# "sum of quantity and pop up for file upload"
# Replace the two sample data rows (row 2 and row 3) with new item data.
# Numeric-looking codes (Article Number, Size, Brand, ItemMRP) are entered
# with a leading apostrophe so Excel keeps them as text, matching how the
# source system exports these reference codes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 -----------------------------------------------------------
$ws.Range("A2").Value = "TROUSER"
$ws.Range("B2").Value = "C1"
$ws.Range("C2").Value = "SC1"
$ws.Range("D2").Value = "'4114"
$ws.Range("F2").Value = "pcs"
$ws.Range("G2").Value = 42
$ws.Range("I2").Value = "L. GREY"
$ws.Range("J2").Value = "'38"
$ws.Range("L2").Value = "'3333"
$ws.Range("R2").Value = "'2198"
$ws.Range("T2").Value = 5

# --- Row 3 -----------------------------------------------------------
$ws.Range("A3").Value = "TROUSER"
$ws.Range("B3").Value = "C1"
$ws.Range("C3").Value = "SC1"
$ws.Range("D3").Value = "'4114"
$ws.Range("F3").Value = "pcs"
$ws.Range("G3").Value = 27
$ws.Range("I3").Value = "L GREY"
$ws.Range("J3").Value = "'36"
$ws.Range("L3").Value = "'1111"
$ws.Range("R3").Value = "'2195.00"
$ws.Range("T3").Value = 1
